$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 5 (G5): replace the shared formula result with a plain literal 0 ---
$ws.Range("G5").Value = 0

# --- Append the new participants (rows 8-14) ---
# Row 8
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Tchomb Joseph Loic"
$ws.Range("C8").Value = 657489972
$ws.Range("D8").Value = 350
$ws.Range("E8").Value = "momo"
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0

# Row 9
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Bayeck Joseph"
$ws.Range("C9").Value = "695919585, 675366970"
$ws.Range("D9").Value = 500
$ws.Range("E9").Value = "cash"
$ws.Range("F9").Value = 200
$ws.Range("G9").Value = 0

# Row 10
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Mabou Fotso Kevin"
$ws.Range("C10").Value = 657936031
$ws.Range("D10").Value = 500
$ws.Range("E10").Value = "cash"
$ws.Range("F10").Value = 200
$ws.Range("G10").Value = 0

# Row 11
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Rikam Yvanol"
$ws.Range("C11").Value = 6553722422
$ws.Range("D11").Value = 300
$ws.Range("E11").Value = "cash"
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0

# Row 12
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Mbame Mbame Martin"
$ws.Range("C12").Value = 697720509
$ws.Range("D12").Value = 300
$ws.Range("E12").Value = "cash"
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0

# Row 13
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Noua Aimerick"
$ws.Range("C13").Value = 674743317
$ws.Range("D13").Value = 300
$ws.Range("E13").Value = "cash"
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0

# Row 14
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Hans Manfred ngalle"
$ws.Range("C14").Value = 650469243
$ws.Range("D14").Value = 1000
$ws.Range("E14").Value = "cash"
$ws.Range("F14").Value = 700
$ws.Range("G14").Value = 0

# --- Update selection to match the author's final cursor position ---
$ws.Range("G10").Select() | Out-Null

# --- Configure page setup (orientation) to match the target layout ---
$ws.PageSetup.Orientation = 1
